$wb = $excel.ActiveWorkbook

# --- Sheet 1 "WithTable": Numbers/Strings/DateTime/Boolean in columns A-D ---
$ws1 = $wb.Worksheets.Item(1)
foreach ($r in @(2,3,4,5)) {
    # Boolean column no longer carries the (redundant) "General" custom style.
    $ws1.Range("D$r").NumberFormat = "General"
}
# Last "Numbers" row becomes a decimal, shown with 2 decimal places.
$ws1.Range("A5").Value = 4.269
$ws1.Range("A5").NumberFormat = "0.00"
[void]$ws1.Range("A5").Select()

# --- Sheet 2 "Tableless": identical layout to sheet 1 ---
$ws2 = $wb.Worksheets.Item(2)
foreach ($r in @(2,3,4,5)) {
    $ws2.Range("D$r").NumberFormat = "General"
}
$ws2.Range("A5").Value = 4.269
$ws2.Range("A5").NumberFormat = "0.00"

# --- Sheet 3 "WithTable_Duplicate": same data, shifted to columns B-G, rows 4-8 ---
$ws3 = $wb.Worksheets.Item(3)
foreach ($r in @(5,6,7,8)) {
    $ws3.Range("E$r").NumberFormat = "General"
}
$ws3.Range("B8").Value = 4.269
$ws3.Range("B8").NumberFormat = "0.00"
[void]$ws3.Range("B8").Select()

# "Tableless" is selected last, making it the active sheet/tab.
[void]$ws2.Range("A5").Select()
